{"js": "// Update the date label and the 25 division-problem cells to the new\n// values from the latest generated output (commit 503736d).\nconst replacements = [\n  [\"2025-01-14 Tuesday\", \"2025-01-15 Wednesday\"],\n  [\"63\u00f73=\", \"40\u00f73=\"],\n  [\"30\u00f76=\", \"92\u00f73=\"],\n  [\"89\u00f75=\", \"15\u00f77=\"],\n  [\"72\u00f74=\", \"20\u00f75=\"],\n  [\"72\u00f76=\", \"68\u00f78=\"],\n  [\"72\u00f77=\", \"46\u00f77=\"],\n  [\"22\u00f76=\", \"57\u00f74=\"],\n  [\"75\u00f75=\", \"61\u00f79=\"],\n  [\"39\u00f76=\", \"18\u00f79=\"],\n  [\"70\u00f73=\", \"56\u00f79=\"],\n  [\"24\u00f72=\", \"13\u00f72=\"],\n  [\"97\u00f74=\", \"91\u00f78=\"],\n  [\"55\u00f78=\", \"62\u00f77=\"],\n  [\"65\u00f77=\", \"95\u00f75=\"],\n  [\"71\u00f75=\", \"34\u00f78=\"],\n  [\"88\u00f79=\", \"28\u00f74=\"],\n  [\"81\u00f73=\", \"86\u00f72=\"],\n  [\"69\u00f78=\", \"40\u00f72=\"],\n  [\"12\u00f78=\", \"98\u00f72=\"],\n  [\"74\u00f78=\", \"84\u00f75=\"],\n  [\"75\u00f79=\", \"36\u00f76=\"],\n  [\"43\u00f77=\", \"80\u00f76=\"],\n  [\"21\u00f74=\", \"45\u00f79=\"],\n  [\"88\u00f77=\", \"68\u00f74=\"],\n  [\"90\u00f76=\", \"42\u00f73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const result of results.items) {\n    result.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date label and the 25 division-problem cells to the new\n# values from the latest generated output (commit 503736d).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-01-14 Tuesday\", \"2025-01-15 Wednesday\"),\n    @(\"63\u00f73=\", \"40\u00f73=\"),\n    @(\"30\u00f76=\", \"92\u00f73=\"),\n    @(\"89\u00f75=\", \"15\u00f77=\"),\n    @(\"72\u00f74=\", \"20\u00f75=\"),\n    @(\"72\u00f76=\", \"68\u00f78=\"),\n    @(\"72\u00f77=\", \"46\u00f77=\"),\n    @(\"22\u00f76=\", \"57\u00f74=\"),\n    @(\"75\u00f75=\", \"61\u00f79=\"),\n    @(\"39\u00f76=\", \"18\u00f79=\"),\n    @(\"70\u00f73=\", \"56\u00f79=\"),\n    @(\"24\u00f72=\", \"13\u00f72=\"),\n    @(\"97\u00f74=\", \"91\u00f78=\"),\n    @(\"55\u00f78=\", \"62\u00f77=\"),\n    @(\"65\u00f77=\", \"95\u00f75=\"),\n    @(\"71\u00f75=\", \"34\u00f78=\"),\n    @(\"88\u00f79=\", \"28\u00f74=\"),\n    @(\"81\u00f73=\", \"86\u00f72=\"),\n    @(\"69\u00f78=\", \"40\u00f72=\"),\n    @(\"12\u00f78=\", \"98\u00f72=\"),\n    @(\"74\u00f78=\", \"84\u00f75=\"),\n    @(\"75\u00f79=\", \"36\u00f76=\"),\n    @(\"43\u00f77=\", \"80\u00f76=\"),\n    @(\"21\u00f74=\", \"45\u00f79=\"),\n    @(\"88\u00f77=\", \"68\u00f74=\"),\n    @(\"90\u00f76=\", \"42\u00f73=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    # 0 = wdFindContinue, 2 = wdReplaceAll\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 0, $false, $newText, 2)\n}\n"}
